# Generate Report for Handoff
# Updates the "Latest Handoff Date(time)" columns for the rows that were
# re-handed-off, on the Overview sheet and on each locale sheet.

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

# Overview sheet: column D holds "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("D$r").Value = "2016-03-25 07:09:59"
}

# zh-cn sheet: column E holds "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "2016-03-25 07:09:54"
}

# de-de sheet: column E holds "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "2016-03-25 07:09:59"
}
